# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 0
    6  = 0
    7  = 8
    8  = 4
    9  = 3
    10 = 1
    11 = 8
    12 = 3
    13 = 3
    14 = 7
    15 = 5
    16 = 0
    17 = 3
    18 = 1
    19 = 1
    20 = 2
    21 = 1
    22 = 0
    23 = 2
    24 = 4
    25 = 2
    26 = 0
    27 = 2
    28 = 0
    29 = 1
    30 = 4
    31 = 1
    32 = 0
    33 = 1
    34 = 2
    35 = 2
    36 = 3
    37 = 5
    38 = 0
    39 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
